{"js": "// Map of old arithmetic expressions -> new ones, as described by the diff.\nconst replacements = [\n  [\"346\u00d73=\", \"145\u00d78=\"],\n  [\"161\u00d74=\", \"638\u00d79=\"],\n  [\"153\u00d73=\", \"800\u00d78=\"],\n  [\"655\u00d76=\", \"612\u00d79=\"],\n  [\"258\u00d72=\", \"267\u00d75=\"],\n  [\"216\u00d75=\", \"713\u00d78=\"],\n  [\"185\u00d79=\", \"375\u00d78=\"],\n  [\"829\u00d78=\", \"874\u00d72=\"],\n  [\"602\u00d77=\", \"188\u00d73=\"],\n  [\"901\u00d75=\", \"458\u00d76=\"],\n  [\"233\u00d75=\", \"129\u00d77=\"],\n  [\"474\u00d75=\", \"341\u00d76=\"],\n  [\"970\u00d72=\", \"257\u00d76=\"],\n  [\"715\u00d75=\", \"142\u00d76=\"],\n  [\"710\u00d78=\", \"566\u00d74=\"],\n  [\"693\u00d75=\", \"904\u00d74=\"],\n  [\"564\u00d79=\", \"744\u00d75=\"],\n  [\"380\u00d72=\", \"750\u00d77=\"],\n  [\"832\u00d76=\", \"450\u00d79=\"],\n  [\"698\u00d72=\", \"498\u00d77=\"],\n  [\"973\u00d75=\", \"723\u00d77=\"],\n  [\"516\u00d73=\", \"978\u00d79=\"],\n  [\"456\u00d74=\", \"417\u00d77=\"],\n  [\"985\u00d72=\", \"984\u00d75=\"],\n  [\"637\u00d72=\", \"412\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 multiplication expressions in the practice table, as\n# described by the diff (each old expression is unique in the document,\n# so a simple Find/Replace per pair is sufficient).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"346\u00d73=\", \"145\u00d78=\"),\n    @(\"161\u00d74=\", \"638\u00d79=\"),\n    @(\"153\u00d73=\", \"800\u00d78=\"),\n    @(\"655\u00d76=\", \"612\u00d79=\"),\n    @(\"258\u00d72=\", \"267\u00d75=\"),\n    @(\"216\u00d75=\", \"713\u00d78=\"),\n    @(\"185\u00d79=\", \"375\u00d78=\"),\n    @(\"829\u00d78=\", \"874\u00d72=\"),\n    @(\"602\u00d77=\", \"188\u00d73=\"),\n    @(\"901\u00d75=\", \"458\u00d76=\"),\n    @(\"233\u00d75=\", \"129\u00d77=\"),\n    @(\"474\u00d75=\", \"341\u00d76=\"),\n    @(\"970\u00d72=\", \"257\u00d76=\"),\n    @(\"715\u00d75=\", \"142\u00d76=\"),\n    @(\"710\u00d78=\", \"566\u00d74=\"),\n    @(\"693\u00d75=\", \"904\u00d74=\"),\n    @(\"564\u00d79=\", \"744\u00d75=\"),\n    @(\"380\u00d72=\", \"750\u00d77=\"),\n    @(\"832\u00d76=\", \"450\u00d79=\"),\n    @(\"698\u00d72=\", \"498\u00d77=\"),\n    @(\"973\u00d75=\", \"723\u00d77=\"),\n    @(\"516\u00d73=\", \"978\u00d79=\"),\n    @(\"456\u00d74=\", \"417\u00d77=\"),\n    @(\"985\u00d72=\", \"984\u00d75=\"),\n    @(\"637\u00d72=\", \"412\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n"}
